$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: update Price (D) / Volume(1h) (E) figures, and
# swap the dogwifhat / InjectiveProtocol rows (44-45) to match the latest
# ranking order pulled from coinranking.com.
#
# Note: several "Price" strings look numeric (e.g. "6.28", "0.470") but must
# stay as literal text (same as the original inlineStr cells), so those are
# written with a leading apostrophe to force Excel's text/quote-prefix
# handling instead of silently re-typing them as numbers (which would also
# strip significant trailing zeros, e.g. "0.470" -> 0.47).
$ws.Range("D2").Value = "60.148.98"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "2.878.70"
$ws.Range("E3").Value = "  -4.34%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'528.43"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("D6").Value = "'128.85"
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "2.877.10"
$ws.Range("E8").Value = "  -4.29%  "
$ws.Range("D9").Value = "'0.470"
$ws.Range("E9").Value = "  -5.18%  "
$ws.Range("D10").Value = "'6.28"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "'0.140"
$ws.Range("E11").Value = "  -5.30%  "
$ws.Range("D12").Value = "'0.425"
$ws.Range("E12").Value = "  -5.12%  "
$ws.Range("D13").Value = "'0.0000209"
$ws.Range("E13").Value = "  -6.11%  "
$ws.Range("D14").Value = "'32.01"
$ws.Range("E14").Value = "  -5.72%  "
$ws.Range("D15").Value = "3.375.82"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").Value = "60.217.89"
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("E17").Value = "  -3.45%  "
$ws.Range("D18").Value = "2.886.81"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").Value = "'6.30"
$ws.Range("E19").Value = "  -5.09%  "
$ws.Range("D20").Value = "'446.27"
$ws.Range("E20").Value = "  -5.52%  "
$ws.Range("D21").Value = "'12.81"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("D22").Value = "'0.622"
$ws.Range("E22").Value = "  -7.59%  "
$ws.Range("D23").Value = "'6.75"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").Value = "'76.34"
$ws.Range("E24").Value = "  -5.06%  "
$ws.Range("D25").Value = "'11.91"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'7.16"
$ws.Range("E29").Value = "  -7.26%  "
$ws.Range("D30").Value = "'1.88"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'24.12"
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("D32").Value = "'1.08"
$ws.Range("E32").Value = "  -5.56%  "
$ws.Range("D33").Value = "'2.21"
$ws.Range("E33").Value = "  -4.20%  "
$ws.Range("D34").Value = "'5.25"
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("D35").Value = "'53.20"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("D36").Value = "'5.55"
$ws.Range("E36").Value = "  -6.06%  "
$ws.Range("D37").Value = "'428.18"
$ws.Range("E37").Value = "  -6.53%  "
$ws.Range("D38").Value = "'0.0767"
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("D39").Value = "'0.0370"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").Value = "2.839.80"
$ws.Range("E40").Value = "  -11.34%  "
$ws.Range("E41").Value = "  -6.62%  "
$ws.Range("D42").Value = "'7.64"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'25.17"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.25"
$ws.Range("E45").Value = "  -6.05%  "
$ws.Range("D46").Value = "'0.234"
$ws.Range("E46").Value = "  -3.52%  "
$ws.Range("D47").Value = "'0.104"
$ws.Range("E47").Value = "  -3.41%  "
$ws.Range("D48").Value = "'1.86"
$ws.Range("E48").Value = "  -6.04%  "
$ws.Range("D49").Value = "'111.13"
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").Value = "0.0₃0463"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'1.19"
$ws.Range("E51").Value = "  -4.93%  "
